$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells whose new value is a plain numeric-looking string (e.g. "18.77", "1.001").
# These must be forced to Text format first, otherwise Excel/COM will silently
# coerce the assignment to a floating point Double and we lose the exact
# original formatting (trailing zeros, decimal-comma grouping, float rounding, etc).
$textFormatUpdates = @{
    'D5' = '239.17'
    'D6' = '1.001'
    'D7' = '0.4824'
    'D8' = '0.2852'
    'D9' = '0.06551'
    'D11' = '0.07475'
    'D12' = '16.73'
    'D13' = '5.107'
    'D14' = '88.09'
    'D15' = '0.6674'
    'D17' = '13.33'
    'D20' = '0.000007591'
    'D21' = '231.36'
    'D22' = '1.002'
    'D23' = '5.287'
    'D24' = '6.230'
    'D25' = '170.07'
    'D26' = '9.348'
    'D27' = '18.77'
    'D28' = '1.967'
    'D29' = '1.402'
    'D30' = '0.1015'
    'D31' = '4.354'
    'D32' = '4.034'
    'D33' = '0.05109'
    'D34' = '1.213'
    'D35' = '0.7589'
    'D36' = '2.713'
    'D37' = '0.01879'
    'D38' = '2.651'
    'D39' = '0.9209'
    'D40' = '2.078'
    'D41' = '107.35'
    'D42' = '0.4304'
    'D44' = '5.706'
    'D45' = '7.431'
    'D46' = '64.65'
    'D47' = '0.1275'
    'D48' = '1.496'
    'D49' = '9.001'
    'D50' = '33.95'
    'D51' = '0.05676'
}

# Remaining cells: values that can never be misread as a pure number
# (percentages, multi-dot "thousand-grouped" numbers, names, URLs) so a
# plain .Value assignment round-trips losslessly as text.
$plainUpdates = @{
    'D2' = '30.689.75'
    'E2' = '  +1.43%  '
    'D3' = '1.896.28'
    'E3' = '  +2.21%  '
    'E4' = '  +0.16%  '
    'E5' = '  +1.26%  '
    'E7' = '  +0.64%  '
    'E8' = '  +1.20%  '
    'E9' = '  +0.77%  '
    'D10' = '1.973.29'
    'E10' = '  +5.96%  '
    'E11' = '  +1.91%  '
    'E12' = '  +2.62%  '
    'E13' = '  -0.69%  '
    'E14' = '  +0.84%  '
    'E15' = '  +2.88%  '
    'D16' = '30.673.13'
    'E16' = '  +1.54%  '
    'B17' = 'Avalanche'
    'C17' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    'E17' = '  +0.36%  '
    'B18' = 'WrappedliquidstakedEther2.0'
    'C18' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D18' = '2.248.43'
    'E18' = '  +7.19%  '
    'E19' = '  +0.19%  '
    'E20' = '  -0.75%  '
    'E21' = '  +3.83%  '
    'B22' = 'BinanceUSD'
    'C22' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    'E22' = '  +0.11%  '
    'B23' = 'Uniswap'
    'C23' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    'E23' = '  -0.29%  '
    'E24' = '  +2.28%  '
    'E25' = '  +3.37%  '
    'E26' = '  +0.84%  '
    'E27' = '  +1.35%  '
    'E28' = '  +2.13%  '
    'E29' = '  -2.21%  '
    'E30' = '  +10.37%  '
    'E31' = '  +2.47%  '
    'E32' = '  +1.57%  '
    'E33' = '  +1.72%  '
    'E34' = '  +5.47%  '
    'E35' = '  +2.33%  '
    'E36' = '  +1.05%  '
    'E37' = '  +3.18%  '
    'E38' = '  +1.41%  '
    'E39' = '  +1.37%  '
    'E40' = '  +0.93%  '
    'E41' = '  +0.70%  '
    'E42' = '  +0.97%  '
    'E43' = '  +0.51%  '
    'E44' = '  -4.59%  '
    'E45' = '  +0.07%  '
    'E46' = '  +0.67%  '
    'E47' = '  -3.55%  '
    'E48' = '  -3.88%  '
    'E49' = '  +2.49%  '
    'E50' = '  -0.88%  '
    'E51' = '  +0.21%  '
}

foreach ($ref in $textFormatUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textFormatUpdates[$ref]
}

foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}
